# Update "paises.xlsx" (Pais sheet) with the latest COVID-19 snapshot.
#
# 1) Refresh the "last updated" timestamp banner in A1.
# 2) Three pairs of countries swapped rank/order (their row in the sheet
#    keeps its statistics lineup, but the displayed country name and,
#    where applicable, the statistics themselves move down/up one slot):
#      - Guyana now sorts ahead of Sierra Leona
#      - Timor Oriental now sorts ahead of Santa Lucia
#      - Islas Malvinas now sorts ahead of Montserrat
# 3) A batch of countries received refreshed case/recovery/death counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh "data as of" banner -----------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 20 de Septiembre de 2020 a las 00:43"

# --- 2. Re-ordered country pairs ---------------------------------------
# Guyana / Sierra Leona (rows 151-152): Guyana moves up with fresh stats,
# Sierra Leona drops down keeping its previous stats.
$ws.Cells.Item(151, 1).Value = "Guyana"
$ws.Cells.Item(151, 2).Value = 2168
$ws.Cells.Item(151, 3).Value = 66
$ws.Cells.Item(151, 4).Value = 1331
$ws.Cells.Item(151, 5).Value = 773
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 2
$ws.Cells.Item(151, 8).Value = 64

$ws.Cells.Item(152, 1).Value = "Sierra Leona"
$ws.Cells.Item(152, 2).Value = 2159
$ws.Cells.Item(152, 3).Value = 6
$ws.Cells.Item(152, 4).Value = 1650
$ws.Cells.Item(152, 5).Value = 437
$ws.Cells.Item(152, 6).Value = 0
$ws.Cells.Item(152, 7).Value = 0
$ws.Cells.Item(152, 8).Value = 72

# Timor Oriental / Santa Lucia (rows 204-205): stats are identical for
# both countries, only the country names trade places.
$ws.Cells.Item(204, 1).Value = "Timor Oriental"
$ws.Cells.Item(205, 1).Value = "Santa Lucia"

# Islas Malvinas / Montserrat (rows 214-215): Islas Malvinas moves up
# with its stats, Montserrat drops down keeping its previous stats.
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(215, 2).Value = 13
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 1

# --- 3. Bulk statistic refresh for unaffected-order countries ----------
# Columns: Row, Casos totales, Nuevos casos, Casos activos, Recuperados,
#          Casos criticos, Muertes hoy, Muertes
$updates = @(
    @(4,   6963319, 37378, 4215954, 2543645, 0, 553, 203720),  # Estados Unidos
    @(8,    762865,  6453,  607837,  123659, 0,  86,  31369),  # Peru
    @(9,    758398,  7927,  627685,  106674, 0, 189,  24039),  # Colombia
    @(13,   622934,  9276,  478077,  132058, 0, 143,  12799),  # Argentina
    @(25,   272308,  1064,  243500,   19342, 0,   2,   9466),  # Alemania
    @(37,   101900,   128,   88666,    7484, 0,  17,   5750),  # Egipto
    @(45,    85152,   808,   74497,    7550, 0,  29,   3105),  # Guatemala
    @(48,    78073,   579,   70495,    6083, 0,  13,   1495),  # Japon
    @(54,    64499,   620,   57299,    6979, 0,   1,    221),  # Barein
    @(58,    57145,   189,   48431,    7619, 0,   1,   1095),  # Nigeria
    @(84,    18819,    86,   13558,    4506, 0,   2,    755),  # Bulgaria
    @(108,    7672,    25,    5914,    1533, 0,   1,    225),  # Zimbabue
    @(115,    5245,    30,    4571,     570, 0,   1,    104)   # Suazilandia
)

foreach ($row in $updates) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}
